# Add filereader to enable file uploads:
# Insert a new "Meta" worksheet at the front of the workbook, holding an
# "Inleiding" (introduction) header + description, used by the CV file
# reader to show an intro blurb. Also nudge a few leftover cell selections
# on the pre-existing sheets (cosmetic cursor-position state saved by Excel).

$wb = $excel.ActiveWorkbook

# --- Insert the new "Meta" sheet before the first existing sheet --------
$firstSheet = $wb.Worksheets.Item(1)
$meta = $wb.Worksheets.Add($firstSheet)
$meta.Name = "Meta"

$meta.Range("A1").Value = "Inleiding"
$meta.Range("A2").Value = "Gemotiveerde en technisch georiënteerde IT-professional met oog voor detail. Via een éénjarig traineeship klaargestoomd voor een carrière als .NET-developer. Gestart bij afdeling Toeslagen en snel doorgegroeid naar medior niveau. Mijn kwaliteiten liggen in het identificeren en doorgronden van problemen, mijn leergierigheid en mijn vermogen om zelfstandig te kunnen werken, maar ook afstemming te zoeken waar dat nodig is."

# Meta becomes the active sheet / tab, cursor parked on F26.
$meta.Activate()
$meta.Range("F26").Select()

# --- Leftover selection nudges on the pre-existing sheets ----------------
$adres = $wb.Worksheets.Item("Adres")
$adres.Range("E17").Select()

$certificaten = $wb.Worksheets.Item("Certificaten")
$certificaten.Range("D20").Select()

$opleidingen = $wb.Worksheets.Item("Opleidingen")
$opleidingen.Range("H23").Select()

# Leave the cursor back on Meta (matches tabSelected on the Meta sheet /
# activeTab="0" on the workbook, and removes tabSelected from Werkervaring,
# which used to be the active tab).
$meta.Activate()
$meta.Range("F26").Select()
